# #5: fund, bonds, otherbonds, antique done
#
# Sheet "基金受益憑證" (fund/beneficiary-certificate sheet, index 5) gains new
# columns: dealer, quantity, face_value (inserted right after owner/name) plus
# the standard trailing block (property_category, category, date,
# legislator_name, legislator_id, source_file, index).
#
# Sheet "具有相當價值之財產" (other valuable property, index 6) gains a new
# "quantity" column plus the same standard trailing block.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet5: 基金受益憑證 (fund)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Existing layout: B=name C=owner D=dealer(face_value co) E=quantity F=face_value G=currency H=total
# Insert 7 new columns starting at D, pushing the old D:H (dealer,quantity,face_value,currency,total) to K:O
$ws5.Columns("D:J").Insert()

# Row 1 - headers
$ws5.Range("B1").Value = "name"
$ws5.Range("C1").Value = "owner"
$ws5.Range("D1").Value = "dealer"
$ws5.Range("E1").Value = "quantity"
$ws5.Range("F1").Value = "face_value"
$ws5.Range("G1").Value = "currency"
$ws5.Range("H1").Value = "total"
$ws5.Range("I1").Value = "property_category"
$ws5.Range("J1").Value = "category"
$ws5.Range("K1").Value = "date"
$ws5.Range("L1").Value = "legislator_name"
$ws5.Range("M1").Value = "legislator_id"
$ws5.Range("N1").Value = "source_file"
$ws5.Range("O1").Value = "index"

# Row 2 - data
$ws5.Range("A2").Value = 64
$ws5.Range("B2").Value = "全球大樂退"
$ws5.Range("C2").Value = "何欣純"
$ws5.Range("D2").Value = "全球人壽"
$ws5.Range("E2").Value = 106.3034
$ws5.Range("F2").Value = 68.33
$ws5.Range("G2").Value = "美金"
$ws5.Range("H2").Value = 212294.31
$ws5.Range("I2").Value = "fund"
$ws5.Range("J2").Value = "normal"
$ws5.Range("K2").NumberFormat = "@"
$ws5.Range("K2").Value = "2012-04-30"
$ws5.Range("L2").Value = "何欣純"
$ws5.Range("M2").Value = 1733
$ws5.Range("N2").Value = "tmp2e891"
$ws5.Range("O2").Value = 64

# ---------------------------------------------------------------------------
# Sheet6: 具有相當價值之財產 (other valuable property / bonds)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Existing layout: B=name C=quantity D=owner E=total
# Insert 7 new columns starting at C, pushing the old C:E (quantity,owner,total) to J:L
$ws6.Columns("C:I").Insert()

# Row 1 - headers
$ws6.Range("B1").Value = "name"
$ws6.Range("C1").Value = "quantity"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "property_category"
$ws6.Range("G1").Value = "category"
$ws6.Range("H1").Value = "date"
$ws6.Range("I1").Value = "legislator_name"
$ws6.Range("J1").Value = "legislator_id"
$ws6.Range("K1").Value = "source_file"
$ws6.Range("L1").Value = "index"

# Row 2 - data
$ws6.Range("A2").Value = 75
$ws6.Range("B2").Value = "國泰富貴保本三福終生壽險"
$ws6.Range("C2").Value = 1
$ws6.Range("D2").Value = "謝俊雄"
$ws6.Range("E2").Value = 331536
$ws6.Range("F2").Value = "otherbonds"
$ws6.Range("G2").Value = "normal"
$ws6.Range("H2").NumberFormat = "@"
$ws6.Range("H2").Value = "2012-04-30"
$ws6.Range("I2").Value = "何欣純"
$ws6.Range("J2").Value = 1733
$ws6.Range("K2").Value = "tmp2e891"
$ws6.Range("L2").Value = 75
